$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Otros" header column (I1), styled like the other headers   ---
# --- (yellow fill) but with only left/right thin borders             ---
$ws.Range("I1").Value = "Otros"
$ws.Range("I1").Interior.Color = 65535
$ws.Range("I1").Borders.Item(7).LineStyle = 1
$ws.Range("I1").Borders.Item(10).LineStyle = 1

# --- Replace row 40 with the new "ASD" product entry ---
$ws.Range("A40").Value = 40
$ws.Range("B40").Value = "ASD"
$ws.Range("C40").Value = 0.77
$ws.Range("D40").Value = 6.01
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = "CONTROL DE ACCESO Y SEGURIDAD"
$ws.Range("G40").Value = "sdf"
$ws.Range("H40").Value = "w"
$ws.Range("I40").Value = "🔒"

# --- Move the active selection, as left by the author ---
$ws.Range("J7").Select() | Out-Null
